$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on Price cells whose new value would otherwise be
# auto-parsed as a number by Excel, so the original textual representation
# (e.g. trailing zeros like "271.90") is preserved, matching the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "43.745.96"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.324.08"
$ws.Range("E3").Value = "  +4.33%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "96.59"
$ws.Range("E5").Value = "  +10.00%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "271.90"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("D10").Value = "45.16"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("D12").Value = "8.12"
$ws.Range("E12").Value = "  +7.29%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "2.672.06"
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("D15").Value = "15.61"
$ws.Range("E15").Value = "  +3.89%  "
$ws.Range("E16").Value = "  +8.32%  "
$ws.Range("D17").Value = "2.332.31"
$ws.Range("E17").Value = "  +3.65%  "
$ws.Range("D18").Value = "43.655.80"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  +5.03%  "
$ws.Range("D20").Value = "6.33"
$ws.Range("E20").Value = "  +6.46%  "
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").Value = "238.12"
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").Value = "2.27"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("D24").Value = "9.51"
$ws.Range("E24").Value = "  +10.10%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  +4.84%  "
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "38.84"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  +8.61%  "
$ws.Range("D32").Value = "172.40"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "0.0896"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "5.49"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").Value = "0.126"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("D36").Value = "0.0360"
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  +9.13%  "
$ws.Range("E41").Value = "  +10.96%  "
$ws.Range("D42").Value = "1.35"
$ws.Range("E42").Value = "  +19.43%  "
$ws.Range("D43").Value = "12.11"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "61.95"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "9.10"
$ws.Range("E45").Value = "  +7.60%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +5.09%  "
$ws.Range("D48").Value = "100.72"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").Value = "2.549.33"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("E51").Value = "  +13.65%  "
